$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, pushing the old row 53 down to row 54.
# The old row 53 (Choclo / Dulce o Americano / Segunda) ends up at row 54
# unchanged, which already matches the target state, so it needs no
# further edits.
$ws.Rows.Item(53).Insert()

# The new (blank) row 53 becomes a copy of what used to be in row 52
# (Choclo / Dulce o Americano / Primera), since the former row 52 is about
# to be overwritten with brand-new data.
$ws.Range("A53").Value2 = 11
$ws.Range("B53").Value2 = "Vega Monumental Concepción"
$ws.Range("C53").Value2 = "Bíobío"
$ws.Range("D53").Value2 = 44292
$ws.Range("E53").Value2 = 8
$ws.Range("F53").Value2 = 100112024
$ws.Range("G53").Value2 = "Choclo"
$ws.Range("H53").Value2 = "Dulce o Americano"
$ws.Range("I53").Value2 = "Primera"
$ws.Range("J53").Value2 = 10000
$ws.Range("K53").Value2 = 150
$ws.Range("L53").Value2 = 200
$ws.Range("M53").Value2 = 175
$ws.Range("N53").Value2 = "$/unidad"
$ws.Range("O53").Value2 = "Región Metropolitana"
$ws.Range("P53").Value2 = 175
$ws.Range("Q53").Value2 = 1
$ws.Range("R53").Value2 = "Hortaliza"

# Row 52 now gets entirely new data (new market record).
$ws.Range("D52").Value2 = 44474
$ws.Range("H52").Value2 = "Choclero"
$ws.Range("I52").Value2 = "Primera"
$ws.Range("J52").Value2 = 100
$ws.Range("K52").Value2 = 39000
$ws.Range("L52").Value2 = 40000
$ws.Range("M52").Value2 = 39500
$ws.Range("N52").Value2 = "$/malla 70 unidades"
$ws.Range("O52").Value2 = "Región de Arica y Parinacota"
$ws.Range("P52").Value2 = 564
$ws.Range("Q52").Value2 = 70
